$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.67479
$ws.Range("C2").Value = 1.03571
$ws.Range("D2").Value = 1.01714
$ws.Range("E2").Value = 1.01258
$ws.Range("F2").Value = 0.99473
$ws.Range("G2").Value = 0.99758
$ws.Range("H2").Value = 1.00349
$ws.Range("I2").Value = 1.01772
$ws.Range("J2").Value = 0.94862

$ws.Range("B3").Value = 0.67138
$ws.Range("C3").Value = 1.0281
$ws.Range("D3").Value = 0.99403
$ws.Range("E3").Value = 0.99618
$ws.Range("F3").Value = 0.9931
$ws.Range("G3").Value = 0.9918
$ws.Range("H3").Value = 0.98875
$ws.Range("I3").Value = 1.00337
$ws.Range("J3").Value = 0.93776

$ws.Range("B4").Value = 0.66865
$ws.Range("C4").Value = 0.99338
$ws.Range("D4").Value = 0.93764
$ws.Range("E4").Value = 0.95246
$ws.Range("F4").Value = 0.91464
$ws.Range("G4").Value = 0.9132
$ws.Range("H4").Value = 0.92706
$ws.Range("I4").Value = 0.96612
$ws.Range("J4").Value = 0.89804

$ws.Range("B5").Value = 0.66737
$ws.Range("C5").Value = 0.9999
$ws.Range("D5").Value = 0.80825
$ws.Range("E5").Value = 0.97752
$ws.Range("F5").Value = 0.8324
$ws.Range("G5").Value = 0.86126
$ws.Range("H5").Value = 0.90572
$ws.Range("I5").Value = 0.9599
$ws.Range("J5").Value = 0.88207

$ws.Range("B6").Value = 0.53369
$ws.Range("C6").Value = 1.05548
$ws.Range("D6").Value = 0.99923
$ws.Range("E6").Value = 1.02379
$ws.Range("F6").Value = 0.98333
$ws.Range("G6").Value = 0.99124
$ws.Range("H6").Value = 1.01316
$ws.Range("I6").Value = 1.01639
$ws.Range("J6").Value = 0.98531

$ws.Range("B7").Value = 0.53045
$ws.Range("C7").Value = 1.04589
$ws.Range("D7").Value = 0.97812
$ws.Range("E7").Value = 1.02121
$ws.Range("F7").Value = 1.02234
$ws.Range("G7").Value = 1.00457
$ws.Range("H7").Value = 1.00145
$ws.Range("I7").Value = 1.00183
$ws.Range("J7").Value = 0.94725

$ws.Range("B8").Value = 0.52726
$ws.Range("C8").Value = 1.03142
$ws.Range("D8").Value = 0.96995
$ws.Range("E8").Value = 1.01705
$ws.Range("F8").Value = 0.99953
$ws.Range("G8").Value = 0.99537
$ws.Range("H8").Value = 0.96205
$ws.Range("I8").Value = 1.01523
$ws.Range("J8").Value = 0.92442

$ws.Range("B9").Value = 0.52563
$ws.Range("C9").Value = 1.03295
$ws.Range("D9").Value = 0.90274
$ws.Range("E9").Value = 1.0117
$ws.Range("F9").Value = 0.98926
$ws.Range("G9").Value = 0.97805
$ws.Range("H9").Value = 0.97813
$ws.Range("I9").Value = 1.02834
$ws.Range("J9").Value = 0.95154

$ws.Range("B10").Value = 1.36129
$ws.Range("C10").Value = 1.02407
$ws.Range("D10").Value = 1.02295
$ws.Range("E10").Value = 1.00005
$ws.Range("F10").Value = 0.98371
$ws.Range("G10").Value = 0.98669
$ws.Range("H10").Value = 0.98841
$ws.Range("I10").Value = 0.99709
$ws.Range("J10").Value = 0.91122

$ws.Range("B11").Value = 1.35567
$ws.Range("C11").Value = 1.02175
$ws.Range("D11").Value = 1.00743
$ws.Range("E11").Value = 0.9796
$ws.Range("F11").Value = 0.9632
$ws.Range("G11").Value = 0.96212
$ws.Range("H11").Value = 0.98184
$ws.Range("I11").Value = 0.99673
$ws.Range("J11").Value = 0.93664

$ws.Range("B12").Value = 1.35293
$ws.Range("C12").Value = 0.98229
$ws.Range("D12").Value = 0.90058
$ws.Range("E12").Value = 0.88685
$ws.Range("F12").Value = 0.78792
$ws.Range("G12").Value = 0.80651
$ws.Range("H12").Value = 0.90208
$ws.Range("I12").Value = 0.91971
$ws.Range("J12").Value = 0.87896

$ws.Range("B13").Value = 1.35223
$ws.Range("C13").Value = 0.98555
$ws.Range("D13").Value = 0.68373
$ws.Range("E13").Value = 0.93797
$ws.Range("F13").Value = 0.64982
$ws.Range("G13").Value = 0.70373
$ws.Range("H13").Value = 0.83365
$ws.Range("I13").Value = 0.90502
$ws.Range("J13").Value = 0.82535

$ws.Range("B14").Value = 0.52134
$ws.Range("C14").Value = 1.02229
$ws.Range("D14").Value = 1.04047
$ws.Range("E14").Value = 1.02048
$ws.Range("F14").Value = 1.04614
$ws.Range("G14").Value = 1.0382
$ws.Range("H14").Value = 1.02118
$ws.Range("I14").Value = 1.07203
$ws.Range("J14").Value = 0.96145

$ws.Range("B15").Value = 0.51862
$ws.Range("C15").Value = 1.00504
$ws.Range("D15").Value = 0.99278
$ws.Range("E15").Value = 0.98372
$ws.Range("F15").Value = 1.00556
$ws.Range("G15").Value = 1.0377
$ws.Range("H15").Value = 0.97883
$ws.Range("I15").Value = 1.02389
$ws.Range("J15").Value = 0.91941

$ws.Range("B16").Value = 0.5159
$ws.Range("C16").Value = 0.93577
$ws.Range("D16").Value = 0.95814
$ws.Range("E16").Value = 0.9735
$ws.Range("F16").Value = 1.01744
$ws.Range("G16").Value = 0.98124
$ws.Range("H16").Value = 0.91349
$ws.Range("I16").Value = 0.97599
$ws.Range("J16").Value = 0.88945

$ws.Range("B17").Value = 0.51436
$ws.Range("C17").Value = 0.96381
$ws.Range("D17").Value = 0.89215
$ws.Range("E17").Value = 1.0041
$ws.Range("F17").Value = 0.88909
$ws.Range("G17").Value = 0.9625
$ws.Range("H17").Value = 0.92565
$ws.Range("I17").Value = 0.94724
$ws.Range("J17").Value = 0.87051

$ws.Range("B18").Value = 1.21295
$ws.Range("C18").Value = 1.03142
$ws.Range("D18").Value = 0.98474
$ws.Range("E18").Value = 0.96907
$ws.Range("F18").Value = 0.97231
$ws.Range("G18").Value = 0.97687
$ws.Range("H18").Value = 0.97776
$ws.Range("I18").Value = 1.00065
$ws.Range("J18").Value = 0.91902

$ws.Range("B19").Value = 1.20767
$ws.Range("C19").Value = 1.02413
$ws.Range("D19").Value = 0.9452
$ws.Range("E19").Value = 0.95298
$ws.Range("F19").Value = 0.89706
$ws.Range("G19").Value = 0.90664
$ws.Range("H19").Value = 0.94543
$ws.Range("I19").Value = 0.9799
$ws.Range("J19").Value = 0.90593

$ws.Range("B20").Value = 1.20367
$ws.Range("C20").Value = 0.97587
$ws.Range("D20").Value = 0.8485
$ws.Range("E20").Value = 0.90311
$ws.Range("F20").Value = 0.82168
$ws.Range("G20").Value = 0.8273
$ws.Range("H20").Value = 0.87429
$ws.Range("I20").Value = 0.92064
$ws.Range("J20").Value = 0.87192

$ws.Range("B21").Value = 1.20212
$ws.Range("C21").Value = 0.98708
$ws.Range("D21").Value = 0.65139
$ws.Range("E21").Value = 0.94958
$ws.Range("F21").Value = 0.69797
$ws.Range("G21").Value = 0.73022
$ws.Range("H21").Value = 0.8409
$ws.Range("I21").Value = 0.90955
$ws.Range("J21").Value = 0.83621

$ws.Range("B22").Value = 0.52722
$ws.Range("C22").Value = 1.05947
$ws.Range("D22").Value = 1.09373
$ws.Range("E22").Value = 1.10227
$ws.Range("F22").Value = 1.05289
$ws.Range("G22").Value = 1.04407
$ws.Range("H22").Value = 1.05265
$ws.Range("I22").Value = 1.05803
$ws.Range("J22").Value = 0.89462

$ws.Range("B23").Value = 0.52302
$ws.Range("C23").Value = 1.03981
$ws.Range("D23").Value = 1.10523
$ws.Range("E23").Value = 1.10037
$ws.Range("F23").Value = 1.17277
$ws.Range("G23").Value = 1.12414
$ws.Range("H23").Value = 1.08943
$ws.Range("I23").Value = 1.06679
$ws.Range("J23").Value = 0.93486

$ws.Range("B24").Value = 0.51967
$ws.Range("C24").Value = 1.00336
$ws.Range("D24").Value = 1.11417
$ws.Range("E24").Value = 1.05742
$ws.Range("F24").Value = 1.011
$ws.Range("G24").Value = 0.98765
$ws.Range("H24").Value = 1.04776
$ws.Range("I24").Value = 1.07181
$ws.Range("J24").Value = 0.88734

$ws.Range("B25").Value = 0.51809
$ws.Range("C25").Value = 1.01549
$ws.Range("D25").Value = 1.06111
$ws.Range("E25").Value = 1.08505
$ws.Range("F25").Value = 0.99439
$ws.Range("G25").Value = 0.99213
$ws.Range("H25").Value = 1.02224
$ws.Range("I25").Value = 1.07546
$ws.Range("J25").Value = 0.89135

$ws.Range("B26").Value = 0.39178
$ws.Range("C26").Value = 1.02938
$ws.Range("D26").Value = 1.06061
$ws.Range("E26").Value = 1.08133
$ws.Range("F26").Value = 1.01959
$ws.Range("G26").Value = 1.02892
$ws.Range("H26").Value = 1.04999
$ws.Range("I26").Value = 1.04205
$ws.Range("J26").Value = 1.09276

$ws.Range("B27").Value = 0.38985
$ws.Range("C27").Value = 1.03197
$ws.Range("D27").Value = 1.06042
$ws.Range("E27").Value = 1.05065
$ws.Range("F27").Value = 1.13919
$ws.Range("G27").Value = 1.14971
$ws.Range("H27").Value = 1.04677
$ws.Range("I27").Value = 1.02903
$ws.Range("J27").Value = 1.0491

$ws.Range("B28").Value = 0.38798
$ws.Range("C28").Value = 1.04562
$ws.Range("D28").Value = 1.06325
$ws.Range("E28").Value = 1.02656
$ws.Range("F28").Value = 1.1214
$ws.Range("G28").Value = 1.11719
$ws.Range("H28").Value = 0.9967
$ws.Range("I28").Value = 1.02715
$ws.Range("J28").Value = 0.99709

$ws.Range("B29").Value = 0.38675
$ws.Range("C29").Value = 1.03254
$ws.Range("D29").Value = 1.03674
$ws.Range("E29").Value = 0.97613
$ws.Range("F29").Value = 1.08842
$ws.Range("G29").Value = 1.13347
$ws.Range("H29").Value = 1.01844
$ws.Range("I29").Value = 1.02834
$ws.Range("J29").Value = 1.02798

